# Auto-generated edit script: updates crypto Price (D) and
# Volume(1h) (E) columns to match the refreshed data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell whose style represents an untouched, plain
# text "Price" cell (column D, row 2 is never style-risky -
# "27.722.29" can never be misread as a pure number) so we can
# restore that exact look after forcing text entry below.
$plainTextStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = "27.722.29"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "1.849.34"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  -1.69%  "

$ws.Range("D5").Value = "'319.59"
$ws.Range("D5").Style = $plainTextStyle
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("E6").Value = "  -1.81%  "

$ws.Range("D7").Value = "'0.4316"
$ws.Range("D7").Style = $plainTextStyle
$ws.Range("E7").Value = "  -3.07%  "

$ws.Range("E8").Value = "  -2.00%  "

$ws.Range("D9").Value = "'0.07357"
$ws.Range("D9").Style = $plainTextStyle
$ws.Range("E9").Value = "  -1.62%  "

$ws.Range("D10").Value = "'0.8784"
$ws.Range("D10").Style = $plainTextStyle
$ws.Range("E10").Value = "  -1.47%  "

$ws.Range("E11").Value = "  -0.64%  "

$ws.Range("D12").Value = "1.854.16"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("E13").Value = "  -0.80%  "

$ws.Range("D14").Value = "'5.449"
$ws.Range("D14").Style = $plainTextStyle
$ws.Range("E14").Value = "  -2.63%  "

$ws.Range("D15").Value = "'0.07145"
$ws.Range("D15").Style = $plainTextStyle
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "'89.05"
$ws.Range("D16").Style = $plainTextStyle
$ws.Range("E16").Value = "  +4.22%  "

$ws.Range("E17").Value = "  -1.97%  "

$ws.Range("D18").Value = "'0.000009004"
$ws.Range("D18").Style = $plainTextStyle
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = $plainTextStyle
$ws.Range("E19").Value = "  -1.86%  "

$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").Value = "27.722.37"
$ws.Range("E21").Value = "  -0.67%  "

$ws.Range("D22").Value = "'5.225"
$ws.Range("D22").Style = $plainTextStyle
$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("E23").Value = "  -2.15%  "

$ws.Range("D24").Value = "2.080.80"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").Value = "'155.34"
$ws.Range("D26").Style = $plainTextStyle
$ws.Range("E26").Value = "  -2.28%  "

$ws.Range("D27").Value = "'18.61"
$ws.Range("D27").Style = $plainTextStyle
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("D28").Value = "'2.193"
$ws.Range("D28").Style = $plainTextStyle
$ws.Range("E28").Value = "  +9.81%  "

$ws.Range("D29").Value = "'5.388"
$ws.Range("D29").Style = $plainTextStyle
$ws.Range("E29").Value = "  -0.77%  "

$ws.Range("D30").Value = "'119.16"
$ws.Range("D30").Style = $plainTextStyle
$ws.Range("E30").Value = "  +0.50%  "

$ws.Range("D31").Value = "'0.08943"
$ws.Range("D31").Style = $plainTextStyle
$ws.Range("E31").Value = "  -1.58%  "

$ws.Range("D32").Value = "'1.233"
$ws.Range("D32").Style = $plainTextStyle
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").Value = "'0.7796"
$ws.Range("D33").Style = $plainTextStyle
$ws.Range("E33").Value = "  -0.76%  "

$ws.Range("D34").Value = "'4.571"
$ws.Range("D34").Style = $plainTextStyle
$ws.Range("E34").Value = "  -1.12%  "

$ws.Range("D35").Value = "'2.920"
$ws.Range("D35").Style = $plainTextStyle
$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("E36").Value = "  -1.90%  "

$ws.Range("E37").Value = "  -1.21%  "

$ws.Range("D38").Value = "'0.05369"
$ws.Range("D38").Style = $plainTextStyle
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("D40").Value = "'7.325"
$ws.Range("D40").Style = $plainTextStyle
$ws.Range("E40").Value = "  +5.41%  "

$ws.Range("D41").Value = "'2.915"
$ws.Range("D41").Style = $plainTextStyle
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").Value = "'0.5148"
$ws.Range("D42").Style = $plainTextStyle

$ws.Range("D43").Value = "'0.1693"
$ws.Range("D43").Style = $plainTextStyle
$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").Value = "'10.75"
$ws.Range("D45").Style = $plainTextStyle
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("D46").Value = "'108.67"
$ws.Range("D46").Style = $plainTextStyle
$ws.Range("E46").Value = "  -3.55%  "

$ws.Range("D47").Value = "'0.4809"
$ws.Range("D47").Style = $plainTextStyle
$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "'0.06478"
$ws.Range("D48").Style = $plainTextStyle
$ws.Range("E48").Value = "  -2.18%  "

$ws.Range("D49").Value = "'1.696"
$ws.Range("D49").Style = $plainTextStyle
$ws.Range("E49").Value = "  -1.95%  "

$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("E51").Value = "  -4.23%  "
